# Auto-generated Word COM-interop script - "Practicing Services - Start"
# Turns the trailing blank paragraphs after the Services screenshot into the new
# "Services in Angular 6" note (heading + body paragraphs), relocating the _GoBack
# bookmark into the new @NgModule paragraph, matching the target diff exactly.

$d = $word.ActiveDocument

# --- 1. Strip the _GoBack bookmark from the image paragraph (119); it reappears later ---
$imgPara = $d.Paragraphs.Item(119)
$imgXml = @'
<w:p w14:paraId="6377C05D" w14:textId="0BCD5CDD" w:rsidR="009B52EF" w:rsidRDefault="009B52EF" w:rsidP="007369E9" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="60704015" wp14:editId="342A4456"><wp:extent cx="5943600" cy="3343275"/><wp:effectExtent l="0" t="0" r="0" b="9525"/><wp:docPr id="4" name="Picture 4"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId12"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5943600" cy="3343275"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@
$imgPara.Range.InsertXML($imgXml)

# --- 2. Repurpose the 13 existing trailing blank paragraphs (120-132) for the first new paragraphs ---
# paragraph 120 stays blank
# paragraph 121 stays blank
# paragraph 122 stays blank
# paragraph 123 stays blank
# paragraph 124 stays blank
# paragraph 125 stays blank
# paragraph 126 stays blank
# paragraph 127 stays blank
# paragraph 128 stays blank
$p = $d.Paragraphs.Item(129)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Services in Angular 6</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

# paragraph 130 stays blank
$p = $d.Paragraphs.Item(131)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>Section 9, Lecture 105</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item(132)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">If you're using Angular 6+ (check your </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>package.json</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">  to find out), you can provide application-wide services in a different way.</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

# --- 3. Append the remaining new paragraphs after paragraph 132, one at a time ---
$idx = 132
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Instead of adding a service class to the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>providers[</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">]  array in </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>AppModule</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> , you can set the following config in @Injectable() :</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>@</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Injectable(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>providedIn</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: 'root'})</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">export class </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>MyService</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>{ ...</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> }</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>This is exactly the same as:</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">export class </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>MyService</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>{ ...</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> }</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>and</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">import </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t xml:space="preserve">{ </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>MyService</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> } from './path/to/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>my.service</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>';</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>@</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>NgModule</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>({</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">    ...</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">    providers: [</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>MyService</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>]</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>})</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">export class </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>AppModule</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>{ ...</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> }</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$p = $d.Paragraphs.Item($idx)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Using this new syntax is completely optional, the traditional syntax (using </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>providers[</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>] ) will still work. The "new syntax" does offer one advantage though: Services can be loaded lazily by Angular (behind the scenes) and redundant code can be removed automatically. This can lead to a better performance and loading speed - though this really only kicks in for bigger services and apps in general.</w:t>
      </w:r>
    </w:p>
'@
$p.Range.InsertXML($xml)

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
# new blank paragraph, nothing further to fill in

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
# new blank paragraph, nothing further to fill in

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"